$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole table is stored as literal text, including the date-/amount-
# looking values, so mark the block as Text before writing: this stops
# Excel re-interpreting "09.11.2022" / "59" / "23.69" as a date or number.
$dataRange = $ws.Range("A2:G4")
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "09.11.2022"
$ws.Range("B2").Value = "59"
$ws.Range("C2").Value = "Бригада"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Материал"
$ws.Range("F2").Value = "Skysawa"
$ws.Range("G2").Value = "Бригада Миши"

$ws.Range("A3").Value = "09.11.2022"
$ws.Range("B3").Value = "23.69"
$ws.Range("C3").Value = "Люди"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Почта"
$ws.Range("F3").Value = "Office"
$ws.Range("G3").Value = ""

$ws.Range("A4").Value = "09.11.2022"
$ws.Range("B4").Value = "26"
$ws.Range("C4").Value = "Люди"
$ws.Range("D4").Value = "Андрей - керовник"
$ws.Range("E4").Value = "Почта"
$ws.Range("F4").Value = "Office"
$ws.Range("G4").Value = ""

# Restore the "General" number format the rest of the sheet uses (matches
# the single-style A1:G4 block in the target file).
$dataRange.NumberFormat = "General"

# A handful of cells are genuinely empty text ("") rather than blank cells
# in the target, not merely absent values. Writing a bare "'" forces a real
# (empty) text cell instead of clearing it, but it also leaves the cell on
# a "quote prefix" style variant - paste-special the formats from an
# untouched plain neighbour on the same row back over it so it lands back
# on the same style the rest of the row uses (a plain .Style= reassignment
# rebinds to the generic "Normal" named style instead, xfId 0).
$ws.Range("D2").Value = "'"
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("D3").Value = "'"
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("G3").Value = "'"
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

$ws.Range("G4").Value = "'"
$ws.Range("F4").Copy()
$ws.Range("G4").PasteSpecial(-4122)

$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.747916666666667)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.747916666666667)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(0.984027777777778)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(0.984027777777778)
